$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.235.55"
$ws.Cells.Item(2, 5).Value = "  +11.60%  "
$ws.Cells.Item(3, 4).Value = "1.819.99"
$ws.Cells.Item(3, 5).Value = "  +8.60%  "
$ws.Cells.Item(4, 5).Value = "  +0.24%  "
$ws.Cells.Item(5, 4).Value = "'229.45"
$ws.Cells.Item(5, 5).Value = "  +4.72%  "
$ws.Cells.Item(6, 4).Value = "'0.574"
$ws.Cells.Item(6, 5).Value = "  +8.65%  "
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 5).Value = "  +0.18%  "
$ws.Cells.Item(8, 4).Value = "'31.47"
$ws.Cells.Item(8, 5).Value = "  +7.69%  "
$ws.Cells.Item(9, 4).Value = "'46.81"
$ws.Cells.Item(9, 5).Value = "  +5.75%  "
$ws.Cells.Item(10, 4).Value = "'0.287"
$ws.Cells.Item(10, 5).Value = "  +8.93%  "
$ws.Cells.Item(11, 4).Value = "'0.0675"
$ws.Cells.Item(11, 5).Value = "  +5.25%  "
$ws.Cells.Item(12, 4).Value = "'0.0930"
$ws.Cells.Item(12, 5).Value = "  +3.10%  "
$ws.Cells.Item(13, 4).Value = "2.085.97"
$ws.Cells.Item(13, 5).Value = "  +8.89%  "
$ws.Cells.Item(14, 4).Value = "1.823.15"
$ws.Cells.Item(14, 5).Value = "  +8.75%  "
$ws.Cells.Item(15, 4).Value = "'0.646"
$ws.Cells.Item(15, 5).Value = "  +7.27%  "
$ws.Cells.Item(16, 4).Value = "34.257.97"
$ws.Cells.Item(16, 5).Value = "  +11.70%  "
$ws.Cells.Item(17, 4).Value = "'10.23"
$ws.Cells.Item(17, 5).Value = "  +1.68%  "
$ws.Cells.Item(18, 4).Value = "'4.28"
$ws.Cells.Item(18, 5).Value = "  +6.81%  "
$ws.Cells.Item(19, 4).Value = "'70.30"
$ws.Cells.Item(19, 5).Value = "  +6.57%  "
$ws.Cells.Item(20, 4).Value = "'257.49"
$ws.Cells.Item(20, 5).Value = "  +5.99%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0754"
$ws.Cells.Item(21, 5).Value = "  +4.64%  "
$ws.Cells.Item(22, 4).Value = "'0.999"
$ws.Cells.Item(22, 5).Value = "  +0.05%  "
$ws.Cells.Item(23, 4).Value = "'10.63"
$ws.Cells.Item(23, 5).Value = "  +6.53%  "
$ws.Cells.Item(24, 4).Value = "'4.33"
$ws.Cells.Item(24, 5).Value = "  +2.24%  "
$ws.Cells.Item(25, 4).Value = "'2.22"
$ws.Cells.Item(25, 5).Value = "  +3.26%  "
$ws.Cells.Item(26, 4).Value = "'159.90"
$ws.Cells.Item(26, 5).Value = "  +0.45%  "
$ws.Cells.Item(27, 4).Value = "'16.71"
$ws.Cells.Item(27, 5).Value = "  +5.77%  "
$ws.Cells.Item(28, 5).Value = "  +5.27%  "
$ws.Cells.Item(29, 4).Value = "'7.13"
$ws.Cells.Item(29, 5).Value = "  +6.78%  "
$ws.Cells.Item(30, 5).Value = "  +0.19%  "
$ws.Cells.Item(31, 4).Value = "'3.87"
$ws.Cells.Item(31, 5).Value = "  +12.18%  "
$ws.Cells.Item(32, 4).Value = "'0.0523"
$ws.Cells.Item(32, 5).Value = "  +6.22%  "
$ws.Cells.Item(33, 5).Value = "  +6.45%  "
$ws.Cells.Item(34, 4).Value = "'3.57"
$ws.Cells.Item(34, 5).Value = "  +8.04%  "
$ws.Cells.Item(35, 4).Value = "1.530.96"
$ws.Cells.Item(35, 5).Value = "  +1.80%  "
$ws.Cells.Item(36, 4).Value = "'1.80"
$ws.Cells.Item(36, 5).Value = "  +1.34%  "
$ws.Cells.Item(37, 5).Value = "  +6.05%  "
$ws.Cells.Item(38, 4).Value = "'0.635"
$ws.Cells.Item(38, 5).Value = "  +6.46%  "
$ws.Cells.Item(39, 5).Value = "  +6.69%  "
$ws.Cells.Item(40, 4).Value = "'83.84"
$ws.Cells.Item(40, 5).Value = "  +0.95%  "
$ws.Cells.Item(41, 4).Value = "'2.79"
$ws.Cells.Item(41, 5).Value = "  +4.38%  "
$ws.Cells.Item(42, 5).Value = "  +3.02%  "
$ws.Cells.Item(43, 4).Value = "'0.908"
$ws.Cells.Item(43, 5).Value = "  +8.39%  "
$ws.Cells.Item(44, 5).Value = "  +4.96%  "
$ws.Cells.Item(45, 4).Value = "'0.0521"
$ws.Cells.Item(45, 5).Value = "  +4.52%  "
$ws.Cells.Item(46, 5).Value = "  +6.25%  "
$ws.Cells.Item(47, 4).Value = "1.979.73"
$ws.Cells.Item(47, 5).Value = "  +9.36%  "
$ws.Cells.Item(48, 4).Value = "'5.83"
$ws.Cells.Item(48, 5).Value = "  +5.35%  "
$ws.Cells.Item(49, 4).Value = "'12.07"
$ws.Cells.Item(49, 5).Value = "  +17.04%  "
$ws.Cells.Item(50, 5).Value = "  +0.11%  "
$ws.Cells.Item(51, 4).Value = "'51.57"
$ws.Cells.Item(51, 5).Value = "  +2.80%  "
